$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 4734.9546  # H58: was 6016.8125
$ws.Cells.Item(58, 10).Value = 10090  # J58: was 23250
$ws.Cells.Item(58, 12).Value = 30270  # L58: was 69750
$ws.Cells.Item(58, 14).Value = -30570  # N58: was -70050

$ws.Cells.Item(107, 8).Value = 413.64706  # H107: was 380.86365
$ws.Cells.Item(107, 9).Value = 127  # I107: was 173.61111
$ws.Cells.Item(107, 10).Value = 5000  # J107: was 1313.5
$ws.Cells.Item(107, 11).Value = 127  # K107: was 173.61111
$ws.Cells.Item(107, 12).Value = 5000  # L107: was 1313.5
$ws.Cells.Item(107, 13).Value = 1793  # M107: was 1746.38889
$ws.Cells.Item(107, 14).Value = -8840  # N107: was -5153.5

$ws.Cells.Item(116, 8).Value = 1882.2727  # H116: was 1887.8125
$ws.Cells.Item(116, 10).Value = 1500  # J116: was 1833.3334
$ws.Cells.Item(116, 12).Value = 1500  # L116: was 1833.3334
$ws.Cells.Item(116, 14).Value = -8384  # N116: was -8717.3334

$ws.Cells.Item(132, 8).Value = 3970913  # H132: was 2697304
$ws.Cells.Item(132, 9).Value = 5104891.5  # I132: was 3761459
$ws.Cells.Item(132, 10).Value = 1987.5  # J132: was 1444.4
$ws.Cells.Item(132, 11).Value = 15314674.5  # K132: was 11284377
$ws.Cells.Item(132, 12).Value = 5962.5  # L132: was 4333.200000000001
$ws.Cells.Item(132, 13).Value = -15312144.5  # M132: was -11281847
$ws.Cells.Item(132, 14).Value = -11022.5  # N132: was -9393.200000000001

$ws.Cells.Item(135, 8).Value = 578.2769  # H135: was 604.6613
$ws.Cells.Item(135, 9).Value = 404.64517  # I135: was 423.54236
$ws.Cells.Item(135, 11).Value = 3641.80653  # K135: was 3811.88124
$ws.Cells.Item(135, 13).Value = -1106.80653  # M135: was -1276.88124

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1238.5416  # H2: was 1273.32
$ws.Cells.Item(2, 9).Value = 931.0909  # I2: was 920.1667
$ws.Cells.Item(2, 10).Value = 1498.6923  # J2: was 1599.3077
$ws.Cells.Item(2, 11).Value = 931.0909  # K2: was 920.1667
$ws.Cells.Item(2, 12).Value = 1498.6923  # L2: was 1599.3077
$ws.Cells.Item(2, 13).Value = -818.0909  # M2: was -807.1667
$ws.Cells.Item(2, 14).Value = -1724.6923  # N2: was -1825.3077

$ws.Cells.Item(32, 8).Value = 698.37  # H32: was 1196.44
$ws.Cells.Item(32, 9).Value = 713.3936  # I32: was 1119.4946
$ws.Cells.Item(32, 10).Value = 463  # J32: was 2218.7144
$ws.Cells.Item(32, 11).Value = 713.3936  # K32: was 1119.4946
$ws.Cells.Item(32, 12).Value = 463  # L32: was 2218.7144
$ws.Cells.Item(32, 13).Value = -426.3936  # M32: was -832.4946
$ws.Cells.Item(32, 14).Value = -1037  # N32: was -2792.7144

$ws.Cells.Item(45, 8).Value = 1408.2941  # H45: was 1268.12
$ws.Cells.Item(45, 9).Value = 1307.75  # I45: was 1098.091
$ws.Cells.Item(45, 10).Value = 1497.6666  # J45: was 1401.7142
$ws.Cells.Item(45, 11).Value = 1307.75  # K45: was 1098.091
$ws.Cells.Item(45, 12).Value = 1497.6666  # L45: was 1401.7142
$ws.Cells.Item(45, 13).Value = -930.75  # M45: was -721.0909999999999
$ws.Cells.Item(45, 14).Value = -2251.6666  # N45: was -2155.7142

$ws.Cells.Item(61, 8).Value = 1311.6595  # H61: was 1469.8667
$ws.Cells.Item(61, 9).Value = 839.55554  # I61: was 928.1818
$ws.Cells.Item(61, 10).Value = 2856.7273  # J61: was 2959.5
$ws.Cells.Item(61, 11).Value = 839.55554  # K61: was 928.1818
$ws.Cells.Item(61, 12).Value = 2856.7273  # L61: was 2959.5
$ws.Cells.Item(61, 13).Value = -627.55554  # M61: was -716.1818
$ws.Cells.Item(61, 14).Value = -3280.7273  # N61: was -3383.5

$ws.Cells.Item(74, 8).Value = 793.2286  # H74: was 703.18604
$ws.Cells.Item(74, 9).Value = 655.09375  # I74: was 585.925
$ws.Cells.Item(74, 11).Value = 655.09375  # K74: was 585.925
$ws.Cells.Item(74, 13).Value = 218.90625  # M74: was 288.075

$ws.Cells.Item(77, 8).Value = 793.2286  # H77: was 703.18604
$ws.Cells.Item(77, 9).Value = 655.09375  # I77: was 585.925
$ws.Cells.Item(77, 11).Value = 3275.46875  # K77: was 2929.625
$ws.Cells.Item(77, 13).Value = 1092.53125  # M77: was 1438.375

$ws.Cells.Item(112, 8).Value = 17277.2  # H112: was 17346.75
$ws.Cells.Item(112, 10).Value = 17277.2  # J112: was 17346.75
$ws.Cells.Item(112, 12).Value = 17277.2  # L112: was 17346.75
$ws.Cells.Item(112, 14).Value = -20231.2  # N112: was -20300.75

$ws.Cells.Item(114, 8).Value = 22075.4  # H114: was 24037.6
$ws.Cells.Item(114, 10).Value = 22075.4  # J114: was 24037.6
$ws.Cells.Item(114, 12).Value = 22075.4  # L114: was 24037.6
$ws.Cells.Item(114, 14).Value = -30753.4  # N114: was -32715.6

$ws.Cells.Item(116, 8).Value = 1238.5416  # H116: was 1273.32
$ws.Cells.Item(116, 9).Value = 931.0909  # I116: was 920.1667
$ws.Cells.Item(116, 10).Value = 1498.6923  # J116: was 1599.3077
$ws.Cells.Item(116, 11).Value = 931.0909  # K116: was 920.1667
$ws.Cells.Item(116, 12).Value = 1498.6923  # L116: was 1599.3077
$ws.Cells.Item(116, 13).Value = 1362.9091  # M116: was 1373.8333
$ws.Cells.Item(116, 14).Value = -6086.6923  # N116: was -6187.3077

$ws.Cells.Item(119, 8).Value = 31499.875  # H119: was 31714.285
$ws.Cells.Item(119, 10).Value = 31499.875  # J119: was 31714.285
$ws.Cells.Item(119, 12).Value = 31499.875  # L119: was 31714.285
$ws.Cells.Item(119, 14).Value = -41175.875  # N119: was -41390.285

$ws.Cells.Item(124, 8).Value = 0  # H124: was 13200
$ws.Cells.Item(124, 10).Value = 0  # J124: was 13200
$ws.Cells.Item(124, 12).Value = 0  # L124: was 13200
$ws.Cells.Item(124, 14).ClearContents()  # N124: was -23020

$ws.Cells.Item(125, 8).Value = 60000  # H125: was 29000
$ws.Cells.Item(125, 10).Value = 60000  # J125: was 29000
$ws.Cells.Item(125, 12).Value = 60000  # L125: was 29000
$ws.Cells.Item(125, 14).Value = -69840  # N125: was -38840

$ws.Cells.Item(136, 8).Value = 1311.6595  # H136: was 1469.8667
$ws.Cells.Item(136, 9).Value = 839.55554  # I136: was 928.1818
$ws.Cells.Item(136, 10).Value = 2856.7273  # J136: was 2959.5
$ws.Cells.Item(136, 11).Value = 2518.66662  # K136: was 2784.5454
$ws.Cells.Item(136, 12).Value = 8570.1819  # L136: was 8878.5
$ws.Cells.Item(136, 13).Value = 31.33338000000003  # M136: was -234.5454
$ws.Cells.Item(136, 14).Value = -13670.1819  # N136: was -13978.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1238.5416  # H3: was 1273.32
$ws.Cells.Item(3, 9).Value = 931.0909  # I3: was 920.1667
$ws.Cells.Item(3, 10).Value = 1498.6923  # J3: was 1599.3077
$ws.Cells.Item(3, 11).Value = 931.0909  # K3: was 920.1667
$ws.Cells.Item(3, 12).Value = 1498.6923  # L3: was 1599.3077
$ws.Cells.Item(3, 13).Value = -817.0909  # M3: was -806.1667
$ws.Cells.Item(3, 14).Value = -1726.6923  # N3: was -1827.3077

$ws.Cells.Item(107, 8).Value = 8649.154  # H107: was 10038.909
$ws.Cells.Item(107, 9).Value = 9244.916999999999  # I107: was 11992
$ws.Cells.Item(107, 10).Value = 1500  # J107: was 1250
$ws.Cells.Item(107, 11).Value = 9244.916999999999  # K107: was 11992
$ws.Cells.Item(107, 12).Value = 1500  # L107: was 1250
$ws.Cells.Item(107, 13).Value = -7324.916999999999  # M107: was -10072
$ws.Cells.Item(107, 14).Value = -5340  # N107: was -5090

$ws.Cells.Item(134, 8).Value = 21121.79  # H134: was 22220.592
$ws.Cells.Item(134, 9).Value = 29447.527  # I134: was 29272.861
$ws.Cells.Item(134, 10).Value = 2388.875  # J134: was 2691.2307
$ws.Cells.Item(134, 11).Value = 88342.58099999999  # K134: was 87818.583
$ws.Cells.Item(134, 12).Value = 7166.625  # L134: was 8073.6921
$ws.Cells.Item(134, 13).Value = -85807.58099999999  # M134: was -85283.583
$ws.Cells.Item(134, 14).Value = -12236.625  # N134: was -13143.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4067294.5  # H31: was 4388495
$ws.Cells.Item(31, 9).Value = 1662.5454  # I31: was 2021.64
$ws.Cells.Item(31, 10).Value = 20838026  # J31: was 12824020
$ws.Cells.Item(31, 11).Value = 1662.5454  # K31: was 2021.64
$ws.Cells.Item(31, 12).Value = 20838026  # L31: was 12824020
$ws.Cells.Item(31, 13).Value = -1367.5454  # M31: was -1726.64
$ws.Cells.Item(31, 14).Value = -20838616  # N31: was -12824610

$ws.Cells.Item(34, 8).Value = 4067294.5  # H34: was 4388495
$ws.Cells.Item(34, 9).Value = 1662.5454  # I34: was 2021.64
$ws.Cells.Item(34, 10).Value = 20838026  # J34: was 12824020
$ws.Cells.Item(34, 11).Value = 1662.5454  # K34: was 2021.64
$ws.Cells.Item(34, 12).Value = 20838026  # L34: was 12824020
$ws.Cells.Item(34, 13).Value = -1460.5454  # M34: was -1819.64
$ws.Cells.Item(34, 14).Value = -20838430  # N34: was -12824424

$ws.Cells.Item(134, 8).Value = 850.2273  # H134: was 933.14923
$ws.Cells.Item(134, 9).Value = 810.7193  # I134: was 908.3090999999999
$ws.Cells.Item(134, 10).Value = 1100.4445  # J134: was 1047
$ws.Cells.Item(134, 11).Value = 2432.1579  # K134: was 2724.9273
$ws.Cells.Item(134, 12).Value = 3301.3335  # L134: was 3141
$ws.Cells.Item(134, 13).Value = 102.8420999999998  # M134: was -189.9272999999998
$ws.Cells.Item(134, 14).Value = -8371.333500000001  # N134: was -8211

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 11437.926  # H99: was 9383.223
$ws.Cells.Item(99, 9).Value = 941.6667  # I99: was 1974.75
$ws.Cells.Item(99, 10).Value = 12749.958  # J99: was 11499.929
$ws.Cells.Item(99, 11).Value = 2825.0001  # K99: was 5924.25
$ws.Cells.Item(99, 12).Value = 38249.874  # L99: was 34499.787
$ws.Cells.Item(99, 13).Value = -579.0001000000002  # M99: was -3678.25
$ws.Cells.Item(99, 14).Value = -42741.874  # N99: was -38991.787

$ws.Cells.Item(121, 8).Value = 969230.4399999999  # H121: was 3626.9268
$ws.Cells.Item(121, 9).Value = 100059  # I121: was 25113.75
$ws.Cells.Item(121, 10).Value = 993374.0600000001  # J121: was 1304.027
$ws.Cells.Item(121, 11).Value = 300177  # K121: was 75341.25
$ws.Cells.Item(121, 12).Value = 2980122.18  # L121: was 3912.081
$ws.Cells.Item(121, 13).Value = -298867  # M121: was -74031.25
$ws.Cells.Item(121, 14).Value = -2982742.18  # N121: was -6532.081

$ws.Cells.Item(122, 8).Value = 1133.8214  # H122: was 1238
$ws.Cells.Item(122, 9).Value = 1092.1538  # I122: was 1158.1666
$ws.Cells.Item(122, 10).Value = 1169.9333  # J122: was 1306.4286
$ws.Cells.Item(122, 11).Value = 9829.3842  # K122: was 10423.4994
$ws.Cells.Item(122, 12).Value = 10529.3997  # L122: was 11757.8574
$ws.Cells.Item(122, 13).Value = -7379.3842  # M122: was -7973.499400000001
$ws.Cells.Item(122, 14).Value = -15429.3997  # N122: was -16657.8574

$ws.Cells.Item(131, 8).Value = 1412860  # H131: was 1427261.6
$ws.Cells.Item(131, 10).Value = 1814985.9  # J131: was 1838847.5
$ws.Cells.Item(131, 12).Value = 5444957.699999999  # L131: was 5516542.5
$ws.Cells.Item(131, 14).Value = -5455037.699999999  # N131: was -5526622.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 25001870  # H113: was 27779744
$ws.Cells.Item(113, 10).Value = 1199  # J113: was 1400
$ws.Cells.Item(113, 12).Value = 1199  # L113: was 1400
$ws.Cells.Item(113, 14).Value = -5539  # N113: was -5740

$ws.Cells.Item(126, 8).Value = 2667.8572  # H126: was 1888
$ws.Cells.Item(126, 9).Value = 3000  # I126: was 1868.5
$ws.Cells.Item(126, 10).Value = 2225  # J126: was 1940
$ws.Cells.Item(126, 11).Value = 9000  # K126: was 5605.5
$ws.Cells.Item(126, 12).Value = 6675  # L126: was 5820
$ws.Cells.Item(126, 13).Value = -6530  # M126: was -3135.5
$ws.Cells.Item(126, 14).Value = -11615  # N126: was -10760

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 19517  # H106: was 20207.777
$ws.Cells.Item(106, 10).Value = 19517  # J106: was 20207.777
$ws.Cells.Item(106, 12).Value = 19517  # L106: was 20207.777
$ws.Cells.Item(106, 14).Value = -22041  # N106: was -22731.777

$ws.Cells.Item(110, 8).Value = 10643.5  # H110: was 10644
$ws.Cells.Item(110, 10).Value = 10643.5  # J110: was 10644
$ws.Cells.Item(110, 12).Value = 10643.5  # L110: was 10644
$ws.Cells.Item(110, 14).Value = -18823.5  # N110: was -18824

$ws.Cells.Item(122, 8).Value = 4738  # H122: was 2884.111
$ws.Cells.Item(122, 9).Value = 7634.6665  # I122: was 3383.2727
$ws.Cells.Item(122, 10).Value = 3000  # J122: was 2540.9375
$ws.Cells.Item(122, 11).Value = 22903.9995  # K122: was 10149.8181
$ws.Cells.Item(122, 12).Value = 9000  # L122: was 7622.8125
$ws.Cells.Item(122, 13).Value = -20453.9995  # M122: was -7699.8181
$ws.Cells.Item(122, 14).Value = -13900  # N122: was -12522.8125

$ws.Cells.Item(132, 8).Value = 7034.3145  # H132: was 7019.2856
$ws.Cells.Item(132, 9).Value = 11038.45  # I132: was 10053.909
$ws.Cells.Item(132, 10).Value = 1695.4667  # J132: was 1883.7693
$ws.Cells.Item(132, 11).Value = 33115.35000000001  # K132: was 30161.727
$ws.Cells.Item(132, 12).Value = 5086.4001  # L132: was 5651.3079
$ws.Cells.Item(132, 13).Value = -30585.35000000001  # M132: was -27631.727
$ws.Cells.Item(132, 14).Value = -10146.4001  # N132: was -10711.3079

$ws.Cells.Item(136, 8).Value = 1981.7037  # H136: was 2465.3438
$ws.Cells.Item(136, 9).Value = 1883.2817  # I136: was 2424.9443
$ws.Cells.Item(136, 10).Value = 2680.5  # J136: was 2683.5
$ws.Cells.Item(136, 11).Value = 5649.8451  # K136: was 7274.8329
$ws.Cells.Item(136, 12).Value = 8041.5  # L136: was 8050.5
$ws.Cells.Item(136, 13).Value = -3099.8451  # M136: was -4724.8329
$ws.Cells.Item(136, 14).Value = -13141.5  # N136: was -13150.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 27565  # H119: was 27848.5
$ws.Cells.Item(119, 10).Value = 27565  # J119: was 27848.5
$ws.Cells.Item(119, 12).Value = 27565  # L119: was 27848.5
$ws.Cells.Item(119, 14).Value = -37241  # N119: was -37524.5

$ws.Cells.Item(132, 8).Value = 843.0484  # H132: was 888.9655
$ws.Cells.Item(132, 9).Value = 806.7593000000001  # I132: was 846.17645
$ws.Cells.Item(132, 10).Value = 1088  # J132: was 1200.7142
$ws.Cells.Item(132, 11).Value = 2420.2779  # K132: was 2538.52935
$ws.Cells.Item(132, 12).Value = 3264  # L132: was 3602.1426
$ws.Cells.Item(132, 13).Value = 109.7221  # M132: was -8.52935000000025
$ws.Cells.Item(132, 14).Value = -8324  # N132: was -8662.142599999999
